$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (Avg Daily Volume) updates
$ws.Range("D26").Value = 652988.1935483871
$ws.Range("E26").Value = 994250.7539854796
$ws.Range("G26").Value = 423523
$ws.Range("H26").Value = 977114
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 704264.3870967742
$ws.Range("K26").Value = 1053785.140530006
$ws.Range("M26").Value = 373513
$ws.Range("N26").Value = 966806.5
$ws.Range("O26").Value = 31
$ws.Range("P26").Value = 857977.2580645161
$ws.Range("Q26").Value = 1218611.033344656
$ws.Range("S26").Value = 581141
$ws.Range("T26").Value = 1181526
$ws.Range("U26").Value = 31
$ws.Range("V26").Value = 892981.1935483871
$ws.Range("W26").Value = 1153712.529117181
$ws.Range("Y26").Value = 548457
$ws.Range("Z26").Value = 1278850.5
$ws.Range("AA26").Value = 31
$ws.Range("AB26").Value = 839033.6129032258
$ws.Range("AC26").Value = 1255112.814201568
$ws.Range("AE26").Value = 453080
$ws.Range("AF26").Value = 1185704
$ws.Range("AG26").Value = 31

# Row 27 (Diff_Vol (Ann - Day)) updates
$ws.Range("D27").Value = 204989.064516129
$ws.Range("J27").Value = 153712.8709677419
$ws.Range("V27").Value = -35003.93548387097
$ws.Range("AB27").Value = 18943.64516129032

# Row 28 (# Obs) updates
$ws.Range("D28").Value = 31
$ws.Range("J28").Value = 31
$ws.Range("P28").Value = 31
$ws.Range("V28").Value = 31
$ws.Range("AB28").Value = 31
